$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (RF)
$ws.Range("B3").Value = 0.901
$ws.Range("C3").Value = 0.893
$ws.Range("D3").Value = 0.108
$ws.Range("E3").Value = 0.329
$ws.Range("F3").Value = 0.238
$ws.Range("G3").Value = 0.973

# Row 4 (NN)
$ws.Range("B4").Value = 0.636
$ws.Range("C4").Value = 0.606
$ws.Range("D4").Value = 0.398
$ws.Range("E4").Value = 0.631
$ws.Range("F4").Value = 0.477
$ws.Range("G4").Value = 0.898

# Row 5 (RNN)
$ws.Range("B5").Value = 0.559
$ws.Range("C5").Value = 0.54
$ws.Range("D5").Value = 0.481
$ws.Range("E5").Value = 0.694
$ws.Range("F5").Value = 0.52
$ws.Range("G5").Value = 0.836
